$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: swap the "high FPS" and ".exe" paragraphs so ".exe" comes
# first and "high FPS" comes second (matching the reordering in the
# diff). Swapping the text content in place keeps paragraph identity
# intact instead of physically moving paragraphs around.
# ---------------------------------------------------------------------
$pFps = $d.Paragraphs(6).Range
$pExe = $d.Paragraphs(7).Range
$fpsText = $pFps.Text.Substring(0, $pFps.Text.Length - 1)
$exeText = $pExe.Text.Substring(0, $pExe.Text.Length - 1)

$rFps = $d.Range($pFps.Start, $pFps.End - 1)
$rFps.Text = $exeText

$pExe2 = $d.Paragraphs(7).Range
$rExe = $d.Range($pExe2.Start, $pExe2.End - 1)
$rExe.Text = $fpsText

# ---------------------------------------------------------------------
# Step 2: replace the last paragraph's text (previously the coal
# shovelling / water pouring story) with the new "visual feedback"
# story, split across two runs just like the source edit (done before
# touching bookmarks, since bookmark add/delete on this engine forces
# a run-coalescing pass over the document).
# ---------------------------------------------------------------------
$pLast = $d.Paragraphs(8).Range
$rLast = $d.Range($pLast.Start, $pLast.End - 1)
$rLast.Text = ""

$firstPart = "As a player, I "
$secondPart = "want visual feedback on my coal and temperature situations so that I am immersed in the experience."

$insPoint = $d.Range($pLast.Start, $pLast.Start)
$insPoint.InsertAfter($firstPart)

$afterFirst = $d.Paragraphs(8).Range
$mid = $afterFirst.Start + $firstPart.Length
$insPoint2 = $d.Range($mid, $mid)
$insPoint2.InsertAfter($secondPart)

# ---------------------------------------------------------------------
# Step 3: move the hidden "_GoBack" bookmark from the end of the last
# paragraph to the end of the paragraph that now reads "high FPS"
# (paragraph 7).
# ---------------------------------------------------------------------
$oldMark = $d.Bookmarks("_GoBack")
$oldMark.Delete()

$destPara = $d.Paragraphs(7).Range
# Anchor on the last character (a non-collapsed range) to sidestep
# degenerate-range quirks, capture it, then delete + retype it so the
# bookmark itself collapses to zero width sitting after all the text.
$lastChar = $d.Range($destPara.End - 2, $destPara.End - 1)
$savedChar = $lastChar.Text
$d.Bookmarks.Add("_GoBack", $lastChar) | Out-Null
$newMark = $d.Bookmarks("_GoBack")
$newMark.Range.Text = ""
$newMark2 = $d.Bookmarks("_GoBack")
$newMark2.Range.InsertBefore($savedChar)
